# cn-#14 rollback unintended testing case change
# Revert the SumProduct sheet back to a simple two-column (B:C) SUMPRODUCT
# example: drop the extra D/E "multiplier" columns and their formulas,
# restore C7's value, and fix the SUMPRODUCT formula/result accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SumProduct")

# Remove the D/E columns' contents (the D multiplier values and the
# E = B*C*D formulas) that were part of the unintended test change.
$ws.Range("D5:E7").ClearContents()

# Restore the original C7 value that existed before the D/E columns
# were introduced.
$ws.Range("C7").Value = 11

# Roll the SUMPRODUCT formula back to only use columns B and C.
$ws.Range("C11").Formula = "=SUMPRODUCT(B5:B7,C5:C7)"

# Restore the selection shown in the saved sheet view.
$ws.Range("D5:E7").Select()
